$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in reservation data for row 2 (mesa ocupada) and row 3 (mesa libre)
$ws.Range("A2").Value = 2
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "123"
$ws.Range("C2").Value = "Ocupada"
$ws.Range("D2").Value = 4

$ws.Range("A3").Value = 3
$ws.Range("C3").Value = "Libre"
$ws.Range("D3").Value = 4
